$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 5
$ws.Range("H5").Value = 290.2
$ws.Range("J5").Value = 150
$ws.Range("L5").Value = 150
$ws.Range("N5").Value = -380

# Row 15
$ws.Range("H15").Value = 185.12
$ws.Range("I15").Value = 185.12
$ws.Range("K15").Value = 555.36
$ws.Range("M15").Value = -386.36

# Row 18
$ws.Range("H18").Value = 202.25
$ws.Range("I18").Value = 202.25
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 202.25
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 81.75
$ws.Range("N18").ClearContents()

# Row 135
$ws.Range("H135").Value = 1389671.5
$ws.Range("I135").Value = 1298.421
$ws.Range("J135").Value = 5158112.5
$ws.Range("K135").Value = 11685.789
$ws.Range("L135").Value = 46423012.5
$ws.Range("M135").Value = -9150.789000000001
$ws.Range("N135").Value = -46428082.5

# Row 137
$ws.Range("H137").Value = 1119.6774
$ws.Range("I137").Value = 1180
$ws.Range("J137").Value = 1010
$ws.Range("K137").Value = 3540
$ws.Range("L137").Value = 3030
$ws.Range("M137").Value = -990
$ws.Range("N137").Value = -8130

# Row 138
$ws.Range("H138").Value = 3663.49
$ws.Range("I138").Value = 2848.2068
$ws.Range("J138").Value = 3996.493
$ws.Range("K138").Value = 8544.6204
$ws.Range("L138").Value = 11989.479
$ws.Range("M138").Value = -3404.6204
$ws.Range("N138").Value = -22269.479

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 5963.3
$ws.Range("I32").Value = 5629.596
$ws.Range("J32").Value = 39000
$ws.Range("K32").Value = 5629.596
$ws.Range("L32").Value = 39000
$ws.Range("M32").Value = -5342.596
$ws.Range("N32").Value = -39574

# Row 61
$ws.Range("H61").Value = 2034.375
$ws.Range("I61").Value = 1952.4286
$ws.Range("J61").Value = 2149.1
$ws.Range("K61").Value = 1952.4286
$ws.Range("L61").Value = 2149.1
$ws.Range("M61").Value = -1740.4286
$ws.Range("N61").Value = -2573.1

# Row 74
$ws.Range("H74").Value = 2018.8846
$ws.Range("I74").Value = 1293.5264
$ws.Range("J74").Value = 3987.7144
$ws.Range("K74").Value = 1293.5264
$ws.Range("L74").Value = 3987.7144
$ws.Range("M74").Value = -419.5264
$ws.Range("N74").Value = -5735.7144

# Row 77
$ws.Range("H77").Value = 2018.8846
$ws.Range("I77").Value = 1293.5264
$ws.Range("J77").Value = 3987.7144
$ws.Range("K77").Value = 6467.632
$ws.Range("L77").Value = 19938.572
$ws.Range("M77").Value = -2099.632
$ws.Range("N77").Value = -28674.572

# Row 92
$ws.Range("H92").Value = 22175
$ws.Range("J92").Value = 22175
$ws.Range("L92").Value = 22175
$ws.Range("N92").Value = -27167

# Row 132
$ws.Range("H132").Value = 1737
$ws.Range("I132").Value = 1490.8206
$ws.Range("J132").Value = 4937.3335
$ws.Range("K132").Value = 4472.4618
$ws.Range("L132").Value = 14812.0005
$ws.Range("M132").Value = -1942.4618
$ws.Range("N132").Value = -19872.0005

# Row 134
$ws.Range("H134").Value = 49078.625
$ws.Range("J134").Value = 53232.715
$ws.Range("L134").Value = 53232.715
$ws.Range("N134").Value = -63372.715

# Row 136
$ws.Range("H136").Value = 2034.375
$ws.Range("I136").Value = 1952.4286
$ws.Range("J136").Value = 2149.1
$ws.Range("K136").Value = 5857.2858
$ws.Range("L136").Value = 6447.299999999999
$ws.Range("M136").Value = -3307.2858
$ws.Range("N136").Value = -11547.3

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 50
$ws.Range("H50").Value = 27825.092
$ws.Range("J50").Value = 27825.092
$ws.Range("L50").Value = 27825.092
$ws.Range("N50").Value = -29075.092

# Row 51
$ws.Range("H51").Value = 29805.105
$ws.Range("J51").Value = 30905.389
$ws.Range("L51").Value = 30905.389
$ws.Range("N51").Value = -32377.389

# Row 59
$ws.Range("H59").Value = 31739.375
$ws.Range("J59").Value = 31739.375
$ws.Range("L59").Value = 31739.375
$ws.Range("N59").Value = -34029.375

# Row 60
$ws.Range("H60").Value = 26026.867
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 26026.867
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 26026.867
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -27048.867

# Row 61
$ws.Range("H61").Value = 29805.105
$ws.Range("J61").Value = 30905.389
$ws.Range("L61").Value = 30905.389
$ws.Range("N61").Value = -31601.389

# Row 68
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498

# Row 71
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488

# Row 74
$ws.Range("H74").Value = 38657
$ws.Range("J74").Value = 38657
$ws.Range("L74").Value = 38657
$ws.Range("N74").Value = -40405

# Row 77
$ws.Range("H77").Value = 38657
$ws.Range("J77").Value = 38657
$ws.Range("L77").Value = 115971
$ws.Range("N77").Value = -124707

# Row 105
$ws.Range("H105").Value = 3101.6667
$ws.Range("I105").Value = 3101.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3101.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1354.6667
$ws.Range("N105").ClearContents()

# Row 122
$ws.Range("H122").Value = 6169
$ws.Range("I122").Value = 1166.6666
$ws.Range("J122").Value = 11171.333
$ws.Range("K122").Value = 3499.9998
$ws.Range("L122").Value = 33513.999
$ws.Range("M122").Value = -1049.9998
$ws.Range("N122").Value = -38413.999

# Row 140
$ws.Range("H140").Value = 76308.14
$ws.Range("J140").Value = 76308.14
$ws.Range("L140").Value = 76308.14
$ws.Range("N140").Value = -86668.14

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 760.8200000000001
$ws.Range("I131").Value = 295.33334
$ws.Range("J131").Value = 842.9647
$ws.Range("K131").Value = 886.0000200000001
$ws.Range("L131").Value = 2528.8941
$ws.Range("M131").Value = 4153.99998
$ws.Range("N131").Value = -12608.8941

# Row 137
$ws.Range("H137").Value = 6681.405
$ws.Range("I137").Value = 2361.9
$ws.Range("J137").Value = 8031.25
$ws.Range("K137").Value = 7085.700000000001
$ws.Range("L137").Value = 24093.75
$ws.Range("M137").Value = -1985.700000000001
$ws.Range("N137").Value = -34293.75

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 113
$ws.Range("H113").Value = 1459
$ws.Range("I113").Value = 1316.1
$ws.Range("J113").Value = 2888
$ws.Range("K113").Value = 1316.1
$ws.Range("L113").Value = 2888
$ws.Range("M113").Value = 853.9000000000001
$ws.Range("N113").Value = -7228

# Row 138
$ws.Range("H138").Value = 47704.2
$ws.Range("J138").Value = 47704.2
$ws.Range("L138").Value = 47704.2
$ws.Range("N138").Value = -57984.2

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 68
$ws.Range("H68").Value = 1630.5333
$ws.Range("J68").Value = 1787.8182
$ws.Range("L68").Value = 1787.8182
$ws.Range("N68").Value = -3285.8182

# Row 71
$ws.Range("H71").Value = 1630.5333
$ws.Range("J71").Value = 1787.8182
$ws.Range("L71").Value = 8939.091
$ws.Range("N71").Value = -16427.091

# Row 127
$ws.Range("H127").Value = 50786.875
$ws.Range("J127").Value = 50786.875
$ws.Range("L127").Value = 50786.875
$ws.Range("N127").Value = -60706.875

# Row 133
$ws.Range("H133").Value = 60581.5
$ws.Range("J133").Value = 60581.5
$ws.Range("L133").Value = 60581.5
$ws.Range("N133").Value = -65641.5

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 96
$ws.Range("H96").Value = 1981.25
$ws.Range("I96").Value = 2212.5
$ws.Range("J96").Value = 1750
$ws.Range("K96").Value = 2212.5
$ws.Range("L96").Value = 1750
$ws.Range("M96").Value = -839.5
$ws.Range("N96").Value = -4496

# Row 137
$ws.Range("H137").Value = 58123.875
$ws.Range("J137").Value = 58123.875
$ws.Range("L137").Value = 58123.875
$ws.Range("N137").Value = -68323.875
